$d = $word.ActiveDocument

# 1. Merge the split runs around the stray "_GoBack" bookmark in the
#    "Beim Verlieren..." bullet back into a single run, and drop that bookmark.
$d.Content.Find.Execute(
    "Beim Verlieren des Spiels wird der eingesetzte Betrag gelöscht.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Beim Verlieren des Spiels wird der eingesetzte Betrag gelöscht.", 2) | Out-Null

# 2. Simplify the empty paragraph (that only carried an en-US language mark)
#    right before "User Story 8" down to a bare paragraph.
$rng = $d.Content
$rng.Find.Execute("Bereits gesetzte Beträge gehen beim Abbruch verloren.", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$para = $rng.Paragraphs.First
$emptyPara = $para.Next
$emptyRange = $emptyPara.Range
$emptyRange.Font.Reset()
$emptyRange.ParagraphFormat.Reset()
$emptyRange.LanguageID = 0

# 3. Move the "_GoBack" bookmark so that it now sits right after "User Story 9".
$rng2 = $d.Content
$rng2.Find.Execute("User Story 9", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$endPoint = $d.Range($rng2.End, $rng2.End)
$endPoint.InsertAfter([char]0xE000)
$bmPoint = $d.Range($rng2.End, $rng2.End)
$d.Bookmarks.Add("_GoBack", $bmPoint) | Out-Null
$markerRange = $d.Range($rng2.End, $rng2.End + 1)
$markerRange.Delete()
